$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: invalid email (saif) + valid password
$ws.Range("A3").Value = "saif"
$ws.Range("B3").Value = "Boxer@123"

# Row 4: valid email + invalid password (numeric)
$ws.Range("A4").Value = "saifzane2@gmail.com"
$ws.Range("B4").Value = 123

# Row 5: invalid email (saif) + invalid password (numeric)
$ws.Range("A5").Value = "saif"
$ws.Range("B5").Value = 123

# Row 6: valid email + valid password
$ws.Range("A6").Value = "saifzane2@gmail.com"
$ws.Range("B6").Value = "Boxer@123"
